$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 10-16: replace with newly computed averaged-intensity values for the existing
# sampling schemes (labels/index unchanged; this run recalculated the numbers).
# Rows 17-19: new rows for schemes that were pushed out of the first 7 (HexGrid rows),
# now appended at the bottom after the new Gaussian-Quadrature/Spiral rows were computed.

# Row 10: Gaussian-Quadrature
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.7618749660278448
$ws.Range("D10").Value = 1.342963186989
$ws.Range("E10").Value = 0.9519747158388456
$ws.Range("F10").Value = 1.054547556704468
$ws.Range("G10").Value = 0.7618749660278448
$ws.Range("H10").Value = 1.342963186989
$ws.Range("I10").Value = 0.8963520568295001
$ws.Range("J10").Value = 1.059311935314379
$ws.Range("K10").Value = 0.8977606274731419
$ws.Range("L10").Value = 1.201449708837466
$ws.Range("M10").Value = 0.7618749660278448
$ws.Range("N10").Value = 1.147468951413923
$ws.Range("O10").Value = 1.02784010639004
$ws.Range("P10").Value = 1.02077934425183

# Row 11: Spiral-90deg-10rot-5space
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.6189703961075463
$ws.Range("D11").Value = 1.461931455641121
$ws.Range("E11").Value = 0.9484343659321822
$ws.Range("F11").Value = 1.088055568034605
$ws.Range("G11").Value = 0.6189703961075463
$ws.Range("H11").Value = 1.461931455641121
$ws.Range("I11").Value = 0.8463812374645328
$ws.Range("J11").Value = 1.113133610327921
$ws.Range("K11").Value = 0.8515275816742304
$ws.Range("L11").Value = 1.291269007727922
$ws.Range("M11").Value = 0.6189703961075463
$ws.Range("N11").Value = 1.205182910786652
$ws.Range("O11").Value = 1.029347946428864
$ws.Range("P11").Value = 1.027462902863758

# Row 12: Spiral-90deg-15rot-5space
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.6199162849221839
$ws.Range("D12").Value = 1.459026196222443
$ws.Range("E12").Value = 0.9494616962368416
$ws.Range("F12").Value = 1.08750873796041
$ws.Range("G12").Value = 0.6199162849221839
$ws.Range("H12").Value = 1.459026196222443
$ws.Range("I12").Value = 0.8473729857073778
$ws.Range("J12").Value = 1.113003891487454
$ws.Range("K12").Value = 0.8519856149342236
$ws.Range("L12").Value = 1.289682161071232
$ws.Range("M12").Value = 0.6199162849221839
$ws.Range("N12").Value = 1.204243946229642
$ws.Range("O12").Value = 1.02897822883547
$ws.Range("P12").Value = 1.027244696067771

# Row 13: Spiral-90deg-10rot-3space
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.6191594854192468
$ws.Range("D13").Value = 1.461341888081326
$ws.Range("E13").Value = 0.9485901815926989
$ws.Range("F13").Value = 1.087904925080094
$ws.Range("G13").Value = 0.6191594854192468
$ws.Range("H13").Value = 1.461341888081326
$ws.Range("I13").Value = 0.846580111367041
$ws.Range("J13").Value = 1.113153222149066
$ws.Range("K13").Value = 0.8516239523561192
$ws.Range("L13").Value = 1.291025223311981
$ws.Range("M13").Value = 0.6191594854192468
$ws.Range("N13").Value = 1.204966034837013
$ws.Range("O13").Value = 1.029249120043342
$ws.Range("P13").Value = 1.027422373669697

# Row 14: NoRotation-tilt60deg
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.4634119999999995
$ws.Range("D14").Value = 1.947848000000001
$ws.Range("E14").Value = 0.7659160000000005
$ws.Range("F14").Value = 1.182999999999997
$ws.Range("G14").Value = 0.4634119999999995
$ws.Range("H14").Value = 1.947848000000001
$ws.Range("I14").Value = 0.6829880000000003
$ws.Range("J14").Value = 1.136060000000001
$ws.Range("K14").Value = 0.7768200000000003
$ws.Range("L14").Value = 1.530199999999998
$ws.Range("M14").Value = 0.4634119999999995
$ws.Range("N14").Value = 1.356882000000001
$ws.Range("O14").Value = 1.090043999999999
$ws.Range("P14").Value = 1.0607805

# Row 15: Rotation-NoTilt
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 0.25
$ws.Range("D15").Value = 2.57
$ws.Range("E15").Value = 0.53
$ws.Range("F15").Value = 1.310550000000001
$ws.Range("G15").Value = 0.25
$ws.Range("H15").Value = 2.57
$ws.Range("I15").Value = 0.4608249999999997
$ws.Range("J15").Value = 1.17
$ws.Range("K15").Value = 0.68
$ws.Range("L15").Value = 1.8727625
$ws.Range("M15").Value = 0.25
$ws.Range("N15").Value = 1.55
$ws.Range("O15").Value = 1.1651375
$ws.Range("P15").Value = 1.1055171875

# Row 16: Rotation-60detTilt
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 0.5620568602624005
$ws.Range("D16").Value = 1.9114556407808
$ws.Range("E16").Value = 0.723883000422398
$ws.Range("F16").Value = 1.1762715988992
$ws.Range("G16").Value = 0.5620568602624005
$ws.Range("H16").Value = 1.9114556407808
$ws.Range("I16").Value = 0.6866615920640008
$ws.Range("J16").Value = 1.100670355763201
$ws.Range("K16").Value = 0.8107865309184008
$ws.Range("L16").Value = 1.506247258931197
$ws.Range("M16").Value = 0.5620568602624005
$ws.Range("N16").Value = 1.317669320601599
$ws.Range("O16").Value = 1.093416775091199
$ws.Range("P16").Value = 1.0597541047552

# Row 17: HexGrid-90degTilt5degRes
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9928659333933656
$ws.Range("D17").Value = 0.9947532445060715
$ws.Range("E17").Value = 0.9987131245920812
$ws.Range("F17").Value = 0.9978421068447468
$ws.Range("G17").Value = 0.9928659333933656
$ws.Range("H17").Value = 0.9947532445060715
$ws.Range("I17").Value = 0.9972460587440232
$ws.Range("J17").Value = 0.9946392705666832
$ws.Range("K17").Value = 0.9948985477501557
$ws.Range("L17").Value = 0.9993227691788865
$ws.Range("M17").Value = 0.9928659333933656
$ws.Range("N17").Value = 0.9967331845490763
$ws.Range("O17").Value = 0.9960436023340662
$ws.Range("P17").Value = 0.9962851319470016

# Row 18: HexGrid-90degTilt22p5degRes
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 1.026335940885563
$ws.Range("D18").Value = 0.9497030985109445
$ws.Range("E18").Value = 1.007250453997243
$ws.Range("F18").Value = 0.988909340278569
$ws.Range("G18").Value = 1.026335940885563
$ws.Range("H18").Value = 0.9497030985109445
$ws.Range("I18").Value = 1.013452596606048
$ws.Range("J18").Value = 0.9828398189924936
$ws.Range("K18").Value = 1.008048213936067
$ws.Range("L18").Value = 0.9696829808144757
$ws.Range("M18").Value = 1.026335940885563
$ws.Range("N18").Value = 0.9784767762540939
$ws.Range("O18").Value = 0.9930497084180799
$ws.Range("P18").Value = 0.9932778055026754

# Row 19: HexGrid-60degTilt5degRes
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 1.099346011518488
$ws.Range("D19").Value = 0.8489100163981035
$ws.Range("E19").Value = 1.033553092259564
$ws.Range("F19").Value = 0.9605383807356203
$ws.Range("G19").Value = 1.099346011518488
$ws.Range("H19").Value = 0.8489100163981035
$ws.Range("I19").Value = 1.056740995265274
$ws.Range("J19").Value = 0.9677902454051207
$ws.Range("K19").Value = 1.029714758963691
$ws.Range("L19").Value = 0.8954914636872461
$ws.Range("M19").Value = 1.099346011518488
$ws.Range("N19").Value = 0.9412315543288339
$ws.Range("O19").Value = 0.985586875227944
$ws.Range("P19").Value = 0.9865106205291384

# Apply the same index-column formatting (bold, centered, bordered) used by A2:A16 to the new A17:A19 cells
$ws.Range("A16").Copy()
$ws.Range("A17:A19").PasteSpecial(-4122)
